$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 4) to the bug report table, reusing the same visual
# formatting ("Good" green style) already used for the ucLucrari row above,
# and a new bug description string.
$ws.Range("A4").Value = "ucLucrari"
$ws.Range("B4").Value = "changing the tipe of Lucrare: is setting the first from list"
$ws.Range("C4").Value = $true

# Match styling of row 3 (A3:B3 use the "Good" cell style - green fill,
# dark green font, centered, thin border).
$ws.Range("A4:B4").Font.Color = $ws.Range("A3").Font.Color
$ws.Range("A4:B4").Interior.Color = $ws.Range("A3").Interior.Color
$ws.Range("A4:B4").HorizontalAlignment = -4108

# C4 keeps the plain bordered/centered look (no red/green fill override -
# conditional formatting sqref already covers C4:C7 and will color it
# green automatically since the value is TRUE).
$ws.Range("A4:C4").Borders.LineStyle = 1
$ws.Range("C4").HorizontalAlignment = -4108

# Move the active selection (cosmetic change from the original G11).
$ws.Range("G10").Select() | Out-Null
